# Update odds values for row 4 (match: Oriente Petrolero - Always Ready)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G4").Value  = 2.55
$ws.Range("H4").Value  = 3.4
$ws.Range("I4").Value  = 2.75
$ws.Range("J4").Value  = 3.1
$ws.Range("K4").Value  = 2.2
$ws.Range("L4").Value  = 3.25

$ws.Range("U4").Value  = 1.62
$ws.Range("V4").Value  = 2.2
$ws.Range("W4").Value  = 10
$ws.Range("X4").Value  = 13
$ws.Range("Y4").Value  = 10
$ws.Range("Z4").Value  = 23
$ws.Range("AA4").Value = 19
$ws.Range("AB4").Value = 26
$ws.Range("AD4").Value = 6.5
$ws.Range("AE4").Value = 12

$ws.Range("AH4").Value = 10
$ws.Range("AI4").Value = 15
$ws.Range("AJ4").Value = 10
$ws.Range("AK4").Value = 26
$ws.Range("AL4").Value = 21
$ws.Range("AM4").Value = 26
$ws.Range("AN4").Value = 4.75
$ws.Range("AO4").Value = 13
$ws.Range("AP4").Value = 21
$ws.Range("AQ4").Value = 41

$ws.Range("AX4").Value = 4.75
$ws.Range("AY4").Value = 15
$ws.Range("AZ4").Value = 21
$ws.Range("BA4").Value = 41
$ws.Range("BB4").Value = 67
